{"js": "// Replace the division problems in the table with their updated values.\nconst replacements = [\n  [\"751\u00f77=\", \"841\u00f74=\"],\n  [\"987\u00f78=\", \"743\u00f79=\"],\n  [\"443\u00f78=\", \"373\u00f75=\"],\n  [\"891\u00f75=\", \"777\u00f76=\"],\n  [\"997\u00f74=\", \"792\u00f74=\"],\n  [\"388\u00f76=\", \"825\u00f72=\"],\n  [\"277\u00f77=\", \"545\u00f79=\"],\n  [\"942\u00f78=\", \"903\u00f76=\"],\n  [\"493\u00f72=\", \"291\u00f77=\"],\n  [\"431\u00f79=\", \"134\u00f75=\"],\n  [\"352\u00f79=\", \"590\u00f73=\"],\n  [\"291\u00f78=\", \"641\u00f77=\"],\n  [\"819\u00f78=\", \"495\u00f79=\"],\n  [\"465\u00f78=\", \"269\u00f75=\"],\n  [\"869\u00f74=\", \"903\u00f72=\"],\n  [\"452\u00f78=\", \"183\u00f74=\"],\n  [\"288\u00f79=\", \"895\u00f74=\"],\n  [\"821\u00f79=\", \"475\u00f76=\"],\n  [\"700\u00f76=\", \"554\u00f77=\"],\n  [\"874\u00f78=\", \"735\u00f75=\"],\n  [\"471\u00f78=\", \"689\u00f77=\"],\n  [\"618\u00f75=\", \"101\u00f73=\"],\n  [\"214\u00f75=\", \"690\u00f75=\"],\n  [\"280\u00f79=\", \"546\u00f78=\"],\n  [\"238\u00f76=\", \"243\u00f78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the division problems in the table to their new values.\n$replacements = @(\n    @{ Old = \"751\u00f77=\"; New = \"841\u00f74=\" },\n    @{ Old = \"987\u00f78=\"; New = \"743\u00f79=\" },\n    @{ Old = \"443\u00f78=\"; New = \"373\u00f75=\" },\n    @{ Old = \"891\u00f75=\"; New = \"777\u00f76=\" },\n    @{ Old = \"997\u00f74=\"; New = \"792\u00f74=\" },\n    @{ Old = \"388\u00f76=\"; New = \"825\u00f72=\" },\n    @{ Old = \"277\u00f77=\"; New = \"545\u00f79=\" },\n    @{ Old = \"942\u00f78=\"; New = \"903\u00f76=\" },\n    @{ Old = \"493\u00f72=\"; New = \"291\u00f77=\" },\n    @{ Old = \"431\u00f79=\"; New = \"134\u00f75=\" },\n    @{ Old = \"352\u00f79=\"; New = \"590\u00f73=\" },\n    @{ Old = \"291\u00f78=\"; New = \"641\u00f77=\" },\n    @{ Old = \"819\u00f78=\"; New = \"495\u00f79=\" },\n    @{ Old = \"465\u00f78=\"; New = \"269\u00f75=\" },\n    @{ Old = \"869\u00f74=\"; New = \"903\u00f72=\" },\n    @{ Old = \"452\u00f78=\"; New = \"183\u00f74=\" },\n    @{ Old = \"288\u00f79=\"; New = \"895\u00f74=\" },\n    @{ Old = \"821\u00f79=\"; New = \"475\u00f76=\" },\n    @{ Old = \"700\u00f76=\"; New = \"554\u00f77=\" },\n    @{ Old = \"874\u00f78=\"; New = \"735\u00f75=\" },\n    @{ Old = \"471\u00f78=\"; New = \"689\u00f77=\" },\n    @{ Old = \"618\u00f75=\"; New = \"101\u00f73=\" },\n    @{ Old = \"214\u00f75=\"; New = \"690\u00f75=\" },\n    @{ Old = \"280\u00f79=\"; New = \"546\u00f78=\" },\n    @{ Old = \"238\u00f76=\"; New = \"243\u00f78=\" }\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n}\n"}
